# Add new Policia data for the "2025" sheet (Negociado de Policia victim data)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2025")
$ws.Activate()

# Updated counts for Mujeres (col B) and Hombres (col C) by age group (rows 2-16)
$newValues = @(
    @(2,  280, 214),
    @(3,  516, 129),
    @(4,  147, 50),
    @(5,  67,  14),
    @(6,  44,  6),
    @(7,  41,  9),
    @(8,  23,  8),
    @(9,  17,  8),
    @(10, 13,  2),
    @(11, 12,  0),
    @(12, 8,   0),
    @(13, 5,   2),
    @(14, 2,   0),
    @(15, 15,  1),
    @(16, 73,  27)
)

foreach ($entry in $newValues) {
    $row = $entry[0]
    $mujeres = $entry[1]
    $hombres = $entry[2]
    $ws.Cells.Item($row, 2).Value = $mujeres
    $ws.Cells.Item($row, 3).Value = $hombres
}

# Move the active selection from F16 to E16, as in the authored edit
$ws.Range("E16").Select()
